# This script reproduces a new daily price record that was inserted into the
# "Femacal de La Calera - Pepino ensalada" sheet, just before the existing
# row 288. Inserting the row pushes every following record (old rows
# 288-364) down by one position, and the record that used to be in row 364
# ends up in the new row 365 - exactly matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 288; this shifts rows
# 288:364 down to 289:365 and keeps their contents (and formatting) intact.
$ws.Rows("288:288").Insert()

# Populate the newly inserted row 288 with the new record's data.
$ws.Range("A288").Value = 3
$ws.Range("B288").Value = "Femacal de La Calera"
$ws.Range("C288").Value = "Coquimbo"
$ws.Range("D288").Value = 44722
$ws.Range("E288").Value = 5
$ws.Range("F288").Value = 100112043
$ws.Range("G288").Value = "Pepino ensalada"
$ws.Range("H288").Value = "Sin especificar"
$ws.Range("I288").Value = "Primera"
$ws.Range("J288").Value = 120
$ws.Range("K288").Value = 20000
$ws.Range("L288").Value = 21000
$ws.Range("M288").Value = 20542
$ws.Range("N288").Value = "`$/caja 70 unidades"
$ws.Range("O288").Value = "Región de Arica y Parinacota"
$ws.Range("P288").Value = 293
$ws.Range("Q288").Value = 70
$ws.Range("R288").Value = "Hortaliza"
